# "Java Basics - completed"
# Add the final glossary row ("final") right after the previous last row
# (row 28) of the Term/"Stands for…"/Definition table, and move the
# selection down to where the user would type the next entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New glossary entry: Term = "final" in column A of the new last row.
# (Table1's range already extends far beyond the data, so the table
# auto-grows to include this row - no explicit resize needed.)
$ws.Range("A29").Value = "final"

# Leave the selection where the author's cursor ended up after typing the
# new row (one row below it, back in column A).
$ws.Range("A31").Select()
